$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) from 45203 to 45204 for rows 2 through 45
for ($r = 2; $r -le 45; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
